# Remove the bold formatting from the two single-letter panel labels
# ("A" and "B") on the slide -- they are plain small text boxes named
# "TextBox 30" and "TextBox 31".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $label = $shp.TextFrame.TextRange.Text
        if ($label -eq "A" -or $label -eq "B") {
            $shp.TextFrame.TextRange.Font.Bold = $false
        }
    }
}
